$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column C ("max") is removed entirely; columns D ("prediction") and
# E ("rejection-f") shift left to become the new C and D.
$ws.Range("C1:C4").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# Column B now holds the new prediction scores instead of the old "1" placeholder.
$ws.Range("B2").Value = 1643.945200406717
$ws.Range("B3").Value = 1535.092286422305
$ws.Range("B4").Value = 1723.472007440319

$wb.Save()
